# Add 8 new library rows (85-92) to Sheet1, describing the new L3 whole-body
# and L3 CEPsh (no spikeIns) samples, matching the commit that expands the
# experimental design table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -------------------------------------------------------
# Columns: A=date, B=Request, C=Multiplex id, D=Library id, E=genotype,
#          F=Tissue.Cell-type.details, G=promoter, H=sampleInfo, I=sample.ID

# Write column H (and I, identical text) first, then F, then the rest, so
# that new shared-string entries are appended in the same order as the
# authoring session (H/I descriptive labels first, then the F "*.L3"
# category labels, matching the target workbook's sharedStrings order).

$ws.Range("H85").Value = "L3 WT Treatment (N2), technical replicate #1 (no spikeIns)"
$ws.Range("I85").Value = "L3 WT Treatment (N2), technical replicate #1 (no spikeIns)"

$ws.Range("H86").Value = "L3 WT No treatment (N2), technical replicate #1 (no spikeIns)"
$ws.Range("I86").Value = "L3 WT No treatment (N2), technical replicate #1 (no spikeIns)"

$ws.Range("H87").Value = "L3 WT Treatment (N2), technical replicate #2(no spikeIns)"
$ws.Range("I87").Value = "L3 WT Treatment (N2), technical replicate #2(no spikeIns)"

$ws.Range("H88").Value = "L3 WT No treatment (N2), technical replicate #2(no spikeIns)"
$ws.Range("I88").Value = "L3 WT No treatment (N2), technical replicate #2(no spikeIns)"

$ws.Range("H89").Value = "L3 OS11358 Treatment(no spikeIns)"
$ws.Range("I89").Value = "L3 OS11358 Treatment(no spikeIns)"

$ws.Range("H90").Value = "L3 OS11358 No treatment(no spikeIns)"
$ws.Range("I90").Value = "L3 OS11358 No treatment(no spikeIns)"

$ws.Range("H91").Value = "L3 OS11359 Treatment(no spikeIns)"
$ws.Range("I91").Value = "L3 OS11359 Treatment(no spikeIns)"

$ws.Range("H92").Value = "L3 OS11359 No treatment(no spikeIns)"
$ws.Range("I92").Value = "L3 OS11359 No treatment(no spikeIns)"

$ws.Range("F85").Value = "whole.body.L3"
$ws.Range("F86").Value = "whole.body.L3"
$ws.Range("F87").Value = "whole.body.L3"
$ws.Range("F88").Value = "whole.body.L3"

$ws.Range("F89").Value = "CEPsh.L3"
$ws.Range("F90").Value = "CEPsh.L3"
$ws.Range("F91").Value = "CEPsh.L3"
$ws.Range("F92").Value = "CEPsh.L3"

# Re-used labels (already present in sharedStrings, no new entries created)
$ws.Range("G85").Value = "No_promoter"
$ws.Range("G86").Value = "No_promoter"
$ws.Range("G87").Value = "No_promoter"
$ws.Range("G88").Value = "No_promoter"

$ws.Range("G89").Value = "promoter-CEPsh"
$ws.Range("G90").Value = "promoter-CEPsh"
$ws.Range("G91").Value = "promoter-CEPsh"
$ws.Range("G92").Value = "promoter-CEPsh"

for ($r = 85; $r -le 92; $r++) {
    $ws.Range("E$r").Value = "WT"
    $ws.Range("B$r").Value = 6239
    $ws.Range("C$r").Value = 6165
}

$ws.Range("D85").Value = 66867
$ws.Range("D86").Value = 66868
$ws.Range("D87").Value = 66869
$ws.Range("D88").Value = 66870
$ws.Range("D89").Value = 66871
$ws.Range("D90").Value = 66872
$ws.Range("D91").Value = 66873
$ws.Range("D92").Value = 66874

# Column A: pick up the existing date-style formatting (built-in m/d/yyyy,
# centered) from the header block above, then stamp in the date serial.
$ws.Range("A2").Copy()
$ws.Range("A85:A92").PasteSpecial(-4122)
for ($r = 85; $r -le 92; $r++) {
    $ws.Range("A$r").Value = 43243
}

# --- Column width / view state -------------------------------------------
$ws.Columns.Item(5).ColumnWidth = 51.166666666666664

$ws.Range("F96").Select() | Out-Null
